$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.081940714174692
$ws.Cells.Item(2, 4).Value = 1.083993805722326
$ws.Cells.Item(2, 5).Value = 1.084782500275467
$ws.Cells.Item(2, 6).Value = 1.095481992770079
$ws.Cells.Item(2, 9).Value = 1.065962289605374
$ws.Cells.Item(2, 10).Value = 1.086812809583376
$ws.Cells.Item(2, 11).Value = 1.086657978799703
$ws.Cells.Item(2, 12).Value = 1.087444627943882
$ws.Cells.Item(2, 13).Value = 1.09811668600872
$ws.Cells.Item(2, 14).Value = 1.088356208432156

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.083440155485552
$ws.Cells.Item(3, 4).Value = 1.085053021434649
$ws.Cells.Item(3, 5).Value = 1.086104269216301
$ws.Cells.Item(3, 6).Value = 1.096796010018168
$ws.Cells.Item(3, 9).Value = 1.06648645698643
$ws.Cells.Item(3, 10).Value = 1.087971616283168
$ws.Cells.Item(3, 11).Value = 1.087535537176838
$ws.Cells.Item(3, 12).Value = 1.088584256449371
$ws.Cells.Item(3, 13).Value = 1.099250576518215
$ws.Cells.Item(3, 14).Value = 1.089516660770379

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.08440880018819
$ws.Cells.Item(4, 4).Value = 1.08573724763094
$ws.Cells.Item(4, 5).Value = 1.086957860037292
$ws.Cells.Item(4, 6).Value = 1.097645054008166
$ws.Cells.Item(4, 9).Value = 1.066823604512597
$ws.Cells.Item(4, 10).Value = 1.088719380565419
$ws.Cells.Item(4, 11).Value = 1.088101614103521
$ws.Cells.Item(4, 12).Value = 1.089319438926541
$ws.Cells.Item(4, 13).Value = 1.099982509287112
$ws.Cells.Item(4, 14).Value = 1.09026548696368

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.084815644439123
$ws.Cells.Item(5, 4).Value = 1.086024623287005
$ws.Cells.Item(5, 5).Value = 1.087316313257954
$ws.Cells.Item(5, 6).Value = 1.098001707319338
$ws.Cells.Item(5, 9).Value = 1.066964859658364
$ws.Cells.Item(5, 10).Value = 1.089033253142698
$ws.Cells.Item(5, 11).Value = 1.088339174917411
$ws.Cells.Item(5, 12).Value = 1.089627980871924
$ws.Cells.Item(5, 13).Value = 1.100289795302612
$ws.Cells.Item(5, 14).Value = 1.090579805276021

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.084883933670041
$ws.Cells.Item(6, 4).Value = 1.086072859028846
$ws.Cells.Item(6, 5).Value = 1.087376476071531
$ws.Cells.Item(6, 6).Value = 1.098061574442612
$ws.Cells.Item(6, 9).Value = 1.066988548855258
$ws.Cells.Item(6, 10).Value = 1.089085925311964
$ws.Cells.Item(6, 11).Value = 1.08837903802903
$ws.Cells.Item(6, 12).Value = 1.08967975558602
$ws.Cells.Item(6, 13).Value = 1.100341365587121
$ws.Cells.Item(6, 14).Value = 1.090632552245809

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.084414237920996
$ws.Cells.Item(7, 4).Value = 1.085741088626424
$ws.Cells.Item(7, 5).Value = 1.086962651256344
$ws.Cells.Item(7, 6).Value = 1.097649820739198
$ws.Cells.Item(7, 9).Value = 1.06682549385992
$ws.Cells.Item(7, 10).Value = 1.088723576454483
$ws.Cells.Item(7, 11).Value = 1.088104790040079
$ws.Cells.Item(7, 12).Value = 1.089323563749494
$ws.Cells.Item(7, 13).Value = 1.099986616895862
$ws.Cells.Item(7, 14).Value = 1.090269688811388

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.082447792066821
$ws.Cells.Item(8, 4).Value = 1.084352014651881
$ws.Cells.Item(8, 5).Value = 1.085229550380749
$ws.Cells.Item(8, 6).Value = 1.095926325611112
$ws.Cells.Item(8, 9).Value = 1.066139855181034
$ws.Cells.Item(8, 10).Value = 1.087204864074935
$ws.Cells.Item(8, 11).Value = 1.086954920959266
$ws.Cells.Item(8, 12).Value = 1.087830236730908
$ws.Cells.Item(8, 13).Value = 1.098500259331434
$ws.Cells.Item(8, 14).Value = 1.088748819686062

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.078970110877433
$ws.Cells.Item(9, 4).Value = 1.081895245873316
$ws.Cells.Item(9, 5).Value = 1.08216243804949
$ws.Cells.Item(9, 6).Value = 1.092879763739399
$ws.Cells.Item(9, 9).Value = 1.064916044693965
$ws.Cells.Item(9, 10).Value = 1.08451264707415
$ws.Cells.Item(9, 11).Value = 1.08491503381627
$ws.Cells.Item(9, 12).Value = 1.085181433340128
$ws.Cells.Item(9, 13).Value = 1.095867312548408
$ws.Cells.Item(9, 14).Value = 1.086052779428335

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.076642735038212
$ws.Cells.Item(10, 4).Value = 1.080251065500051
$ws.Cells.Item(10, 5).Value = 1.080108443736594
$ws.Cells.Item(10, 6).Value = 1.090841950224358
$ws.Cells.Item(10, 9).Value = 1.064089484289999
$ws.Cells.Item(10, 10).Value = 1.082706669719427
$ws.Cells.Item(10, 11).Value = 1.083545669417436
$ws.Cells.Item(10, 12).Value = 1.083403512514609
$ws.Cells.Item(10, 13).Value = 1.094102406065205
$ws.Cells.Item(10, 14).Value = 1.084244237378619

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.075632731633793
$ws.Cells.Item(11, 4).Value = 1.079537557366446
$ws.Cells.Item(11, 5).Value = 1.079216757286759
$ws.Cells.Item(11, 6).Value = 1.089957868272329
$ws.Cells.Item(11, 9).Value = 1.063728999339432
$ws.Cells.Item(11, 10).Value = 1.081921929402595
$ws.Cells.Item(11, 11).Value = 1.082950424708443
$ws.Cells.Item(11, 12).Value = 1.082630711043497
$ws.Cells.Item(11, 13).Value = 1.093335828115106
$ws.Cells.Item(11, 14).Value = 1.083458382640526

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.075257226347873
$ws.Cells.Item(12, 4).Value = 1.079272288186446
$ws.Cells.Item(12, 5).Value = 1.078885192879964
$ws.Cells.Item(12, 6).Value = 1.089629219205375
$ws.Cells.Item(12, 9).Value = 1.063594708186667
$ws.Cells.Item(12, 10).Value = 1.081630022922065
$ws.Cells.Item(12, 11).Value = 1.082728973740129
$ws.Cells.Item(12, 12).Value = 1.082343207673774
$ws.Cells.Item(12, 13).Value = 1.093050725555971
$ws.Cells.Item(12, 14).Value = 1.083166061619312

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.075337789271886
$ws.Cells.Item(13, 4).Value = 1.079329200290103
$ws.Cells.Item(13, 5).Value = 1.078956330641597
$ws.Cells.Item(13, 6).Value = 1.089699727499372
$ws.Cells.Item(13, 9).Value = 1.063623531845282
$ws.Cells.Item(13, 10).Value = 1.081692656930875
$ws.Cells.Item(13, 11).Value = 1.082776491668129
$ws.Cells.Item(13, 12).Value = 1.082404898625928
$ws.Cells.Item(13, 13).Value = 1.093111897481369
$ws.Cells.Item(13, 14).Value = 1.083228784575599

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.075601699343112
$ws.Cells.Item(14, 4).Value = 1.079515635068507
$ws.Cells.Item(14, 5).Value = 1.079189357302795
$ws.Cells.Item(14, 6).Value = 1.089930707416378
$ws.Cells.Item(14, 9).Value = 1.063717906787593
$ws.Cells.Item(14, 10).Value = 1.081897808911478
$ws.Cells.Item(14, 11).Value = 1.082932126682523
$ws.Cells.Item(14, 12).Value = 1.082606955157558
$ws.Cells.Item(14, 13).Value = 1.093312268858633
$ws.Cells.Item(14, 14).Value = 1.083434227895546

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.075764257043974
$ws.Cells.Item(15, 4).Value = 1.079630471669822
$ws.Cells.Item(15, 5).Value = 1.079332885794385
$ws.Cells.Item(15, 6).Value = 1.090072986883674
$ws.Cells.Item(15, 9).Value = 1.063776002387474
$ws.Cells.Item(15, 10).Value = 1.08202415407848
$ws.Cells.Item(15, 11).Value = 1.083027971949014
$ws.Cells.Item(15, 12).Value = 1.082731388941621
$ws.Cells.Item(15, 13).Value = 1.093435676194733
$ws.Cells.Item(15, 14).Value = 1.083560752487176

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.076709717765156
$ws.Cells.Item(16, 4).Value = 1.080298385237195
$ws.Cells.Item(16, 5).Value = 1.08016757304621
$ws.Cells.Item(16, 6).Value = 1.090900587520855
$ws.Cells.Item(16, 9).Value = 1.064113353853365
$ws.Cells.Item(16, 10).Value = 1.082758691994848
$ws.Cells.Item(16, 11).Value = 1.083585124954203
$ws.Cells.Item(16, 12).Value = 1.083454738028767
$ws.Cells.Item(16, 13).Value = 1.094153231001411
$ws.Cells.Item(16, 14).Value = 1.084296333531639

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.077302175514243
$ws.Cells.Item(17, 4).Value = 1.080716926649085
$ws.Cells.Item(17, 5).Value = 1.08069053081716
$ws.Cells.Item(17, 6).Value = 1.091419260761987
$ws.Cells.Item(17, 9).Value = 1.064324272588703
$ws.Cells.Item(17, 10).Value = 1.083218709455203
$ws.Cells.Item(17, 11).Value = 1.083933992856363
$ws.Cells.Item(17, 12).Value = 1.083907680868035
$ws.Cells.Item(17, 13).Value = 1.094602697250344
$ws.Cells.Item(17, 14).Value = 1.084757004269552

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.077647531154464
$ws.Cells.Item(18, 4).Value = 1.080960903906524
$ws.Cells.Item(18, 5).Value = 1.080995342492248
$ws.Cells.Item(18, 6).Value = 1.091721630995571
$ws.Cells.Item(18, 9).Value = 1.064447049339369
$ws.Cells.Item(18, 10).Value = 1.083486765963488
$ws.Cells.Item(18, 11).Value = 1.084137259932513
$ws.Cells.Item(18, 12).Value = 1.084171590624206
$ws.Cells.Item(18, 13).Value = 1.094864635875104
$ws.Cells.Item(18, 14).Value = 1.085025441448825

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.077765252327466
$ws.Cells.Item(19, 4).Value = 1.081044068408191
$ws.Cells.Item(19, 5).Value = 1.081099238189249
$ws.Cells.Item(19, 6).Value = 1.091824703977772
$ws.Cells.Item(19, 9).Value = 1.064488870981978
$ws.Cells.Item(19, 10).Value = 1.083578121764528
$ws.Cells.Item(19, 11).Value = 1.084206531301025
$ws.Cells.Item(19, 12).Value = 1.084261529055565
$ws.Cells.Item(19, 13).Value = 1.094953911797241
$ws.Cells.Item(19, 14).Value = 1.085116926985581

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.077238632705026
$ws.Cells.Item(20, 4).Value = 1.080672036758069
$ws.Cells.Item(20, 5).Value = 1.080634445286841
$ws.Cells.Item(20, 6).Value = 1.091363628951911
$ws.Cells.Item(20, 9).Value = 1.064301668741972
$ws.Cells.Item(20, 10).Value = 1.083169381283935
$ws.Cells.Item(20, 11).Value = 1.083896585595885
$ws.Cells.Item(20, 12).Value = 1.083859113830625
$ws.Cells.Item(20, 13).Value = 1.094554497347541
$ws.Cells.Item(20, 14).Value = 1.084707606046622

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.075523994005629
$ws.Cells.Item(21, 4).Value = 1.079460741352427
$ws.Cells.Item(21, 5).Value = 1.079120746548866
$ws.Cells.Item(21, 6).Value = 1.089862696867391
$ws.Cells.Item(21, 9).Value = 1.063690126541244
$ws.Cells.Item(21, 10).Value = 1.081837408391455
$ws.Cells.Item(21, 11).Value = 1.082886305767149
$ws.Cells.Item(21, 12).Value = 1.082547467039757
$ws.Cells.Item(21, 13).Value = 1.093253274510721
$ws.Cells.Item(21, 14).Value = 1.083373741599856

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.074443931550244
$ws.Cells.Item(22, 4).Value = 1.078697757924298
$ws.Cells.Item(22, 5).Value = 1.07816698153633
$ws.Cells.Item(22, 6).Value = 1.088917483879819
$ws.Cells.Item(22, 9).Value = 1.063303361752596
$ws.Cells.Item(22, 10).Value = 1.080997516386045
$ws.Cells.Item(22, 11).Value = 1.082249071936874
$ws.Cells.Item(22, 12).Value = 1.081720172573562
$ws.Cells.Item(22, 13).Value = 1.092433049976217
$ws.Cells.Item(22, 14).Value = 1.082532656851462

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.075016685915817
$ws.Cells.Item(23, 4).Value = 1.079102363697396
$ws.Cells.Item(23, 5).Value = 1.078672786644742
$ws.Cells.Item(23, 6).Value = 1.089418705343736
$ws.Cells.Item(23, 9).Value = 1.063508608871821
$ws.Cells.Item(23, 10).Value = 1.081442991606461
$ws.Cells.Item(23, 11).Value = 1.082587075789171
$ws.Cells.Item(23, 12).Value = 1.082158986873586
$ws.Cells.Item(23, 13).Value = 1.092868067214033
$ws.Cells.Item(23, 14).Value = 1.082978764697787

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.077267345629715
$ws.Cells.Item(24, 4).Value = 1.080692321034001
$ws.Cells.Item(24, 5).Value = 1.080659788606145
$ws.Cells.Item(24, 6).Value = 1.091388767076463
$ws.Cells.Item(24, 9).Value = 1.064311883215551
$ws.Cells.Item(24, 10).Value = 1.083191671376838
$ws.Cells.Item(24, 11).Value = 1.083913489011869
$ws.Cells.Item(24, 12).Value = 1.083881060062455
$ws.Cells.Item(24, 13).Value = 1.094576277511863
$ws.Cells.Item(24, 14).Value = 1.084729927794014

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.079870711787707
$ws.Cells.Item(25, 4).Value = 1.082531477198644
$ws.Cells.Item(25, 5).Value = 1.082956961481275
$ws.Cells.Item(25, 6).Value = 1.093668540412909
$ws.Cells.Item(25, 9).Value = 1.065234299395527
$ws.Cells.Item(25, 10).Value = 1.085210590228878
$ws.Cells.Item(25, 11).Value = 1.085444039028834
$ws.Cells.Item(25, 12).Value = 1.085868309234677
$ws.Cells.Item(25, 13).Value = 1.096549661371605
$ws.Cells.Item(25, 14).Value = 1.086751713742398

Write-Output "Updated vm_pu values for rows 2-25"